# Update the "taxonomic_coverage" sheet with kingdom/phylum/class/order/
# family/genus/species detail for the two existing rows (chinook, steelhead),
# then leave that sheet selected/active (matching the author's final UI
# state) instead of the "funding" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("taxonomic_coverage")

# Row 3 (steelhead) first, then row 2 (chinook) -- matches the order new
# shared strings were minted in the target workbook.
$ws.Range("B3").Value = "steelhead trout "
$ws.Range("C3").Value = "Animalia"
$ws.Range("D3").Value = "Chordata"
$ws.Range("E3").Value = "Actinopterygii"
$ws.Range("F3").Value = "Salmoniformes"
$ws.Range("G3").Value = "Salmonidae"
$ws.Range("H3").Value = "Oncorhynchus"
$ws.Range("I3").Value = "mykiss"

$ws.Range("B2").Value = "chinook salmon"
$ws.Range("C2").Value = "Animalia"
$ws.Range("D2").Value = "Craniata"
$ws.Range("E2").Value = "Actinopterygii"
$ws.Range("F2").Value = "Salmoniformes"
$ws.Range("G2").Value = "Salmonidae"
$ws.Range("H2").Value = "Oncorhynchus"
$ws.Range("I2").Value = "Oncorhynchus tshawytscha"

# Make taxonomic_coverage the active/selected sheet (was "funding"), with
# E15 as the selected cell.
$ws.Select()
$ws.Range("E15").Select()
